$d = $word.ActiveDocument

# --- Edit 1: "Functionalities" heading -> "Functionalities of Game" ---
$rng = $d.Content
$rng.Find.Execute("Functionalities", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rng.Collapse(0)
$rng.InsertAfter(" of Game")
$rng.Bold = 1

# --- Edit 2: "Reward system" heading -> "Reward system/Trophy room:" ---
$rng2 = $d.Content
$rng2.Find.Execute("Reward system", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rng2.Collapse(0)
$rng2.InsertAfter("/Trophy room:")
$rng2.Bold = 1

# --- Edit 3: extend the "reward system purpose" paragraph with new sentences
#     about the trophy room, with a _GoBack bookmark inserted mid-sentence ---
$rng3 = $d.Content
$rng3.Find.Execute("choosing correct options.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rng3.Collapse(0)
$rng3.InsertAfter(" When a player receive reward, it will get added in the trophy room")
$rng3.Collapse(0)
$rng3.InsertAfter(". A player")
$rng3.Collapse(0)
$rng3.InsertAfter(" can visit")
$rng3.Collapse(0)
$rng3.InsertAfter(" the trophy room and check the number of trophies they have achieved.")

# Re-find the mid-paragraph boundary between "the trophy room " and "and check"
# so the _GoBack bookmark can be placed precisely there (a collapsed range at
# the very end of a paragraph's last run gets snapped to the whole paragraph,
# so we insert all the text first, then bookmark a now-interior position).
$rngBookmark = $d.Content
$rngBookmark.Find.Execute("the trophy room ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rngBookmark.Collapse(0)
$d.Bookmarks.Add("_GoBack", $rngBookmark)

# --- Edit 4: "Multiple levels" heading -> "Multiple levels:" (new bold run) ---
$rng4 = $d.Content
$rng4.Find.Execute("Multiple levels", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rng4.Collapse(0)
$rng4.InsertAfter(":")
$rng4.Bold = 1
